$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.937
$ws.Range("B4").Value = 5.942
$ws.Range("A7").Value = -20.523
$ws.Range("C10").Value = -12.859
$ws.Range("B12").Value = 4.939
$ws.Range("C13").Value = -13.097
$ws.Range("A16").Value = -22.109
$ws.Range("B18").Value = 5.231
$ws.Range("B19").Value = 7.754999999999998
$ws.Range("B20").Value = 6.867999999999999
$ws.Range("D25").Value = -7.831
$ws.Range("A28").Value = -21.522
$ws.Range("A29").Value = -21.496
$ws.Range("C30").Value = -11.769
$ws.Range("B31").Value = 6.395999999999999
$ws.Range("A32").Value = -22.002
$ws.Range("D34").Value = -7.904000000000001
$ws.Range("D39").Value = -7.692
$ws.Range("A40").Value = -21.659
$ws.Range("B40").Value = 5.88
$ws.Range("C40").Value = -11.584
$ws.Range("B42").Value = 6.854000000000001
$ws.Range("C44").Value = -11.682
$ws.Range("B47").Value = 5.781000000000001
$ws.Range("B48").Value = 5.523999999999999
$ws.Range("A52").Value = -21.514
$ws.Range("A57").Value = -22.042
$ws.Range("D61").Value = -8.157999999999998
$ws.Range("B63").Value = 5.252
$ws.Range("B64").Value = 5.935
$ws.Range("D64").Value = -7.906999999999999
$ws.Range("A66").Value = -21.482
$ws.Range("B76").Value = 6.462000000000001
$ws.Range("D78").Value = -8.343999999999999
$ws.Range("B81").Value = 5.979000000000001
$ws.Range("D83").Value = -8.1
$ws.Range("B89").Value = 4.81
$ws.Range("C89").Value = -13.71
$ws.Range("C91").Value = -12.353
$ws.Range("D92").Value = -7.117
$ws.Range("B94").Value = 5.887
$ws.Range("D98").Value = -7.468000000000001
$ws.Range("A100").Value = -22.124
